# Auto-generated edit script replicating the diff changes
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 200
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 64.166664
$ws.Range("I38").Value = 64.166664
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 192.499992
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = 179.500008
$ws.Range("H92").Value = 539
$ws.Range("I92").Value = 545.5
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 545.5
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 702.5
$ws.Range("N92").Value = -2996
$ws.Range("H107").Value = 1111.6818
$ws.Range("I107").Value = 944.2353000000001
$ws.Range("J107").Value = 1681
$ws.Range("K107").Value = 944.2353000000001
$ws.Range("L107").Value = 1681
$ws.Range("M107").Value = 975.7646999999999
$ws.Range("N107").Value = -5521
$ws.Range("H112").Value = 3585311.5
$ws.Range("J112").Value = 3704788.8
$ws.Range("L112").Value = 11114366.4
$ws.Range("N112").Value = -11116582.4
$ws.Range("H129").Value = 1696.862
$ws.Range("J129").Value = 1869.6923
$ws.Range("L129").Value = 5609.0769
$ws.Range("N129").Value = -15609.0769
$ws.Range("H132").Value = 2459.5898
$ws.Range("I132").Value = 2567.3333
$ws.Range("J132").Value = 1166.6666
$ws.Range("K132").Value = 7701.999899999999
$ws.Range("L132").Value = 3499.9998
$ws.Range("M132").Value = -5171.999899999999
$ws.Range("N132").Value = -8559.9998

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9067.348
$ws.Range("I32").Value = 6681.738
$ws.Range("J32").Value = 22296.637
$ws.Range("K32").Value = 6681.738
$ws.Range("L32").Value = 22296.637
$ws.Range("M32").Value = -6394.738
$ws.Range("N32").Value = -22870.637
$ws.Range("H45").Value = 2663.2812
$ws.Range("I45").Value = 2395.318
$ws.Range("K45").Value = 2395.318
$ws.Range("M45").Value = -2018.318
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H122").Value = 1894.4572
$ws.Range("I122").Value = 1773.3
$ws.Range("K122").Value = 5319.9
$ws.Range("M122").Value = -2869.9
$ws.Range("H132").Value = 14038.884
$ws.Range("I132").Value = 2047.5883
$ws.Range("J132").Value = 59339.332
$ws.Range("K132").Value = 6142.7649
$ws.Range("L132").Value = 178017.996
$ws.Range("M132").Value = -3612.7649
$ws.Range("N132").Value = -183077.996
$ws.Range("H135").Value = 45886.285
$ws.Range("J135").Value = 45886.285
$ws.Range("L135").Value = 45886.285
$ws.Range("N135").Value = -56026.285

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H36").Value = 670
$ws.Range("I36").Value = 670
$ws.Range("K36").Value = 670
$ws.Range("M36").Value = -136
$ws.Range("H43").Value = 91558
$ws.Range("J43").Value = 91558
$ws.Range("L43").Value = 91558
$ws.Range("N43").Value = -91920
$ws.Range("H54").Value = 8667.666999999999
$ws.Range("I54").Value = 3001.7144
$ws.Range("K54").Value = 3001.7144
$ws.Range("M54").Value = -2517.7144
$ws.Range("H94").Value = 659.6667
$ws.Range("I94").Value = 529.86664
$ws.Range("J94").Value = 789.4666999999999
$ws.Range("K94").Value = 529.86664
$ws.Range("L94").Value = 789.4666999999999
$ws.Range("M94").Value = -78.86663999999996
$ws.Range("N94").Value = -1691.4667
$ws.Range("H99").Value = 1686
$ws.Range("I99").Value = 2166.6667
$ws.Range("J99").Value = 1325.5
$ws.Range("K99").Value = 2166.6667
$ws.Range("L99").Value = 1325.5
$ws.Range("M99").Value = -668.6667000000002
$ws.Range("N99").Value = -4321.5
$ws.Range("H105").Value = 1390891.1
$ws.Range("I105").Value = 1261.1538
$ws.Range("J105").Value = 2176334.2
$ws.Range("K105").Value = 1261.1538
$ws.Range("L105").Value = 2176334.2
$ws.Range("M105").Value = 485.8462
$ws.Range("N105").Value = -2179828.2
$ws.Range("H109").Value = 43685
$ws.Range("J109").Value = 43685
$ws.Range("L109").Value = 43685
$ws.Range("N109").Value = -46459

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3744.45
$ws.Range("I31").Value = 1828.6666
$ws.Range("J31").Value = 6085.963
$ws.Range("K31").Value = 1828.6666
$ws.Range("L31").Value = 6085.963
$ws.Range("M31").Value = -1533.6666
$ws.Range("N31").Value = -6675.963
$ws.Range("H34").Value = 3744.45
$ws.Range("I34").Value = 1828.6666
$ws.Range("J34").Value = 6085.963
$ws.Range("K34").Value = 1828.6666
$ws.Range("L34").Value = 6085.963
$ws.Range("M34").Value = -1626.6666
$ws.Range("N34").Value = -6489.963
$ws.Range("H43").Value = 47000
$ws.Range("J43").Value = 47000
$ws.Range("L43").Value = 47000
$ws.Range("N43").Value = -47368
$ws.Range("H58").Value = 20864
$ws.Range("I58").Value = 1745.4166
$ws.Range("K58").Value = 1745.4166
$ws.Range("M58").Value = -1542.4166
$ws.Range("H86").Value = 11996.546
$ws.Range("I86").Value = 1635.6666
$ws.Range("J86").Value = 15881.875
$ws.Range("K86").Value = 1635.6666
$ws.Range("L86").Value = 15881.875
$ws.Range("M86").Value = -512.6666
$ws.Range("N86").Value = -18127.875
$ws.Range("H89").Value = 11996.546
$ws.Range("I89").Value = 1635.6666
$ws.Range("J89").Value = 15881.875
$ws.Range("K89").Value = 8178.333000000001
$ws.Range("L89").Value = 79409.375
$ws.Range("M89").Value = -2562.333000000001
$ws.Range("N89").Value = -90641.375
$ws.Range("H101").Value = 47000
$ws.Range("J101").Value = 47000
$ws.Range("L101").Value = 47000
$ws.Range("N101").Value = -53490
$ws.Range("H105").Value = 3424
$ws.Range("I105").Value = 4006.6667
$ws.Range("J105").Value = 2550
$ws.Range("K105").Value = 4006.6667
$ws.Range("L105").Value = 2550
$ws.Range("M105").Value = -2259.6667
$ws.Range("N105").Value = -6044
$ws.Range("H109").Value = 141451570
$ws.Range("J109").Value = 141451570
$ws.Range("L109").Value = 141451570
$ws.Range("N109").Value = -141453650
$ws.Range("H134").Value = 1248.579
$ws.Range("I134").Value = 1013.05884
$ws.Range("J134").Value = 3250.5
$ws.Range("K134").Value = 3039.17652
$ws.Range("L134").Value = 9751.5
$ws.Range("M134").Value = -504.17652
$ws.Range("N134").Value = -14821.5
$ws.Range("H136").Value = 20864
$ws.Range("I136").Value = 1745.4166
$ws.Range("K136").Value = 5236.2498
$ws.Range("M136").Value = -2686.2498

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 249.5
$ws.Range("I19").Value = 99
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 297
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = -123
$ws.Range("N19").Value = -1548
$ws.Range("H68").Value = 20676.4
$ws.Range("I68").Value = 689.5
$ws.Range("J68").Value = 34001
$ws.Range("K68").Value = 2068.5
$ws.Range("L68").Value = 102003
$ws.Range("M68").Value = -1257.5
$ws.Range("N68").Value = -103625
$ws.Range("H71").Value = 20676.4
$ws.Range("I71").Value = 689.5
$ws.Range("J71").Value = 34001
$ws.Range("K71").Value = 6205.5
$ws.Range("L71").Value = 306009
$ws.Range("M71").Value = -2149.5
$ws.Range("N71").Value = -314121
$ws.Range("H74").Value = 9953.909
$ws.Range("I74").Value = 9799
$ws.Range("K74").Value = 29397
$ws.Range("M74").Value = -28336
$ws.Range("H77").Value = 9953.909
$ws.Range("I77").Value = 9799
$ws.Range("K77").Value = 88191
$ws.Range("M77").Value = -82887
$ws.Range("H120").Value = 12838
$ws.Range("I120").Value = 5676
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 17028
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -12190
$ws.Range("N120").Value = -69676
$ws.Range("H131").Value = 734.5106
$ws.Range("I131").Value = 338.125
$ws.Range("J131").Value = 771.3837
$ws.Range("K131").Value = 1014.375
$ws.Range("L131").Value = 2314.1511
$ws.Range("M131").Value = 4025.625
$ws.Range("N131").Value = -12394.1511

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 3000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
$ws.Range("H97").Value = 4125
$ws.Range("I97").Value = 3833.3333
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 3833.3333
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -3337.3333
$ws.Range("N97").Value = -5992
$ws.Range("H100").Value = 40400
$ws.Range("J100").Value = 40400
$ws.Range("L100").Value = 40400
$ws.Range("N100").Value = -42564
$ws.Range("H102").Value = 2898.389
$ws.Range("I102").Value = 2571.4375
$ws.Range("K102").Value = 2571.4375
$ws.Range("M102").Value = -949.4375
$ws.Range("H113").Value = 6374.1113
$ws.Range("I113").Value = 10827.223
$ws.Range("J113").Value = 1921
$ws.Range("K113").Value = 10827.223
$ws.Range("L113").Value = 1921
$ws.Range("M113").Value = -8657.223
$ws.Range("N113").Value = -6261

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 148.66667
$ws.Range("I55").Value = 139.08333
$ws.Range("J55").Value = 158.25
$ws.Range("K55").Value = 139.08333
$ws.Range("L55").Value = 158.25
$ws.Range("M55").Value = 33.91667000000001
$ws.Range("N55").Value = -504.25
$ws.Range("H132").Value = 229896.03
$ws.Range("I132").Value = 327821.5
$ws.Range("K132").Value = 983464.5
$ws.Range("M132").Value = -980934.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1037.5333
$ws.Range("I132").Value = 729.08
$ws.Range("J132").Value = 2579.8
$ws.Range("K132").Value = 2187.24
$ws.Range("L132").Value = 7739.400000000001
$ws.Range("M132").Value = 342.7599999999998
$ws.Range("N132").Value = -12799.4
$ws.Range("H136").Value = 38235612
$ws.Range("I136").Value = 54331320
$ws.Range("J136").Value = 8300.625
$ws.Range("K136").Value = 162993960
$ws.Range("L136").Value = 24901.875
$ws.Range("M136").Value = -162991410
$ws.Range("N136").Value = -30001.875
